$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the Server IP for Fileservice (row 5) and FTP (row 9) from the
# old value 20.107.94.83 to the new value 108.143.141.90
$ws.Range("E5").Value = "108.143.141.90"
$ws.Range("E9").Value = "108.143.141.90"

# Widen column E slightly to fit the new (longer) Server IP value
$ws.Columns.Item(5).ColumnWidth = 12.83

# Move the active selection to C8 (as seen when the workbook was saved)
$ws.Range("C8").Select() | Out-Null
